# Scheduled market-data refresh: update Universalis price snapshots and
# recomputed profit figures on each crafting job sheet (Jenova_Profits).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2141.484
$ws.Range("I17").Value = 1300
$ws.Range("J17").Value = 2169.5334
$ws.Range("K17").Value = 3900
$ws.Range("L17").Value = 6508.600199999999
$ws.Range("M17").Value = -3732
$ws.Range("N17").Value = -6844.600199999999
$ws.Range("H41").Value = 21742128
$ws.Range("I41").Value = 515.8333
$ws.Range("J41").Value = 45460252
$ws.Range("K41").Value = 515.8333
$ws.Range("L41").Value = 45460252
$ws.Range("M41").Value = -75.83330000000001
$ws.Range("N41").Value = -45461132
$ws.Range("H86").Value = 2774320.8
$ws.Range("I86").Value = 4430.4287
$ws.Range("K86").Value = 4430.4287
$ws.Range("M86").Value = -3307.4287
$ws.Range("H89").Value = 2774320.8
$ws.Range("I89").Value = 4430.4287
$ws.Range("K89").Value = 22152.1435
$ws.Range("M89").Value = -16536.1435
$ws.Range("H92").Value = 185.85715
$ws.Range("I92").Value = 141.91667
$ws.Range("K92").Value = 141.91667
$ws.Range("M92").Value = 1106.08333
$ws.Range("H100").Value = 2633.3333
$ws.Range("I100").Value = 900
$ws.Range("J100").Value = 3500
$ws.Range("K100").Value = 900
$ws.Range("L100").Value = 3500
$ws.Range("M100").Value = -359
$ws.Range("N100").Value = -4582
$ws.Range("H107").Value = 34807.9
$ws.Range("I107").Value = 36043.895
$ws.Range("K107").Value = 36043.895
$ws.Range("M107").Value = -34123.895
$ws.Range("H112").Value = 3601.1853
$ws.Range("I112").Value = 1979
$ws.Range("K112").Value = 5937
$ws.Range("M112").Value = -4829
$ws.Range("H116").Value = 18458.5
$ws.Range("I116").Value = 6732.4
$ws.Range("K116").Value = 6732.4
$ws.Range("M116").Value = -3290.4
$ws.Range("H132").Value = 4617.875
$ws.Range("I132").Value = 4824
$ws.Range("K132").Value = 14472
$ws.Range("M132").Value = -11942

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4087.3594
$ws.Range("I32").Value = 3596.6826
$ws.Range("K32").Value = 3596.6826
$ws.Range("M32").Value = -3309.6826
$ws.Range("H88").Value = 2396.1
$ws.Range("J88").Value = 2436.7144
$ws.Range("L88").Value = 2436.7144
$ws.Range("N88").Value = -3248.7144
$ws.Range("H91").Value = 2396.1
$ws.Range("J91").Value = 2436.7144
$ws.Range("L91").Value = 2436.7144
$ws.Range("N91").Value = -5244.7144
$ws.Range("H132").Value = 4155.1523
$ws.Range("J132").Value = 5930.7856
$ws.Range("L132").Value = 17792.3568
$ws.Range("N132").Value = -22852.3568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 84302.664
$ws.Range("J59").Value = 84302.664
$ws.Range("L59").Value = 84302.664
$ws.Range("N59").Value = -85996.664
$ws.Range("H99").Value = 4093.75
$ws.Range("I99").Value = 2121.3333
$ws.Range("K99").Value = 2121.3333
$ws.Range("M99").Value = -623.3332999999998
$ws.Range("H105").Value = 1400.2
$ws.Range("I105").Value = 1109.8
$ws.Range("K105").Value = 1109.8
$ws.Range("M105").Value = 637.2
$ws.Range("H107").Value = 715621.9
$ws.Range("I107").Value = 1834.3334
$ws.Range("K107").Value = 1834.3334
$ws.Range("M107").Value = 85.66660000000002
$ws.Range("H134").Value = 31896.914
$ws.Range("I134").Value = 2250.7827
$ws.Range("K134").Value = 6752.348100000001
$ws.Range("M134").Value = -4217.348100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 349.25
$ws.Range("I107").Value = 384.85715
$ws.Range("K107").Value = 384.85715
$ws.Range("M107").Value = 1535.14285
$ws.Range("H134").Value = 872554.8
$ws.Range("I134").Value = 558318.5600000001
$ws.Range("K134").Value = 1674955.68
$ws.Range("M134").Value = -1672420.68
$ws.Range("H139").Value = 73750
$ws.Range("J139").Value = 73750
$ws.Range("L139").Value = 73750
$ws.Range("N139").Value = -84030

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 527096.5
$ws.Range("I92").Value = 769760.25
$ws.Range("K92").Value = 2309280.75
$ws.Range("M92").Value = -2308032.75
$ws.Range("H113").Value = 1684898.6
$ws.Range("I113").Value = 9259834
$ws.Range("J113").Value = 1579.5555
$ws.Range("K113").Value = 27779502
$ws.Range("L113").Value = 4738.666499999999
$ws.Range("M113").Value = -27777332
$ws.Range("N113").Value = -9078.666499999999
$ws.Range("H131").Value = 13416913
$ws.Range("I131").Value = 33433992
$ws.Range("J131").Value = 72194.266
$ws.Range("K131").Value = 100301976
$ws.Range("L131").Value = 216582.798
$ws.Range("M131").Value = -100296936
$ws.Range("N131").Value = -226662.798

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 410883.62
$ws.Range("I132").Value = 592350.5
$ws.Range("K132").Value = 1777051.5
$ws.Range("M132").Value = -1774521.5
$ws.Range("H137").Value = 50000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 50000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 50000
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -60200

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3124.5
$ws.Range("I82").Value = 3124.5
$ws.Range("K82").Value = 3124.5
$ws.Range("M82").Value = -2763.5
$ws.Range("H85").Value = 3124.5
$ws.Range("I85").Value = 3124.5
$ws.Range("K85").Value = 3124.5
$ws.Range("M85").Value = -1876.5
$ws.Range("H122").Value = 590582.7
$ws.Range("I122").Value = 2353.818
$ws.Range("K122").Value = 7061.454000000001
$ws.Range("M122").Value = -4611.454000000001
$ws.Range("H132").Value = 2877.08
$ws.Range("I132").Value = 2046.35
$ws.Range("J132").Value = 6200
$ws.Range("K132").Value = 6139.049999999999
$ws.Range("L132").Value = 18600
$ws.Range("M132").Value = -3609.049999999999
$ws.Range("N132").Value = -23660

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1655.6364
$ws.Range("I81").Value = 1141.5333
$ws.Range("K81").Value = 2283.0666
$ws.Range("M81").Value = -1222.0666
$ws.Range("H84").Value = 1655.6364
$ws.Range("I84").Value = 1141.5333
$ws.Range("K84").Value = 11415.333
$ws.Range("M84").Value = -6111.333000000001
$ws.Range("H122").Value = 50005656
$ws.Range("I122").Value = 66671948
$ws.Range("K122").Value = 200015844
$ws.Range("M122").Value = -200013394
$ws.Range("H132").Value = 78576.42999999999
$ws.Range("I132").Value = 10598.8
$ws.Range("J132").Value = 116341.78
$ws.Range("K132").Value = 31796.4
$ws.Range("L132").Value = 349025.34
$ws.Range("M132").Value = -29266.4
$ws.Range("N132").Value = -354085.34
$ws.Range("H138").Value = 77018.336
$ws.Range("J138").Value = 77018.336
$ws.Range("L138").Value = 77018.336
$ws.Range("N138").Value = -87298.336
